# Weekly update: insert this week's new price record for Achicoria
# (Vega Central Mapocho de Santiago) at the top of the data block (row 16)
# and push the older rows down by one, matching the upstream weekly dump.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 16:40 down to 17:41, creating a blank row 16.
$ws.Rows.Item(16).Insert()

# Populate the newly inserted row 16 with this week's record.
$ws.Cells.Item(16, 1).Value = 9
$ws.Cells.Item(16, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(16, 3).Value = "Metropolitana"
$ws.Cells.Item(16, 4).Value = 44967
$ws.Cells.Item(16, 5).Value = 13
$ws.Cells.Item(16, 6).Value = 100112010
$ws.Cells.Item(16, 7).Value = "Achicoria"
$ws.Cells.Item(16, 8).Value = "Sin especificar"
$ws.Cells.Item(16, 9).Value = "Primera"
$ws.Cells.Item(16, 10).Value = 70
$ws.Cells.Item(16, 11).Value = 7000
$ws.Cells.Item(16, 12).Value = 7000
$ws.Cells.Item(16, 13).Value = 7000
$ws.Cells.Item(16, 14).Value = "$/caja 16 unidades"
$ws.Cells.Item(16, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(16, 16).Value = 438
$ws.Cells.Item(16, 17).Value = 16
$ws.Cells.Item(16, 18).Value = "Hortaliza"
